# Auto-generated edit script: "Add data for 2024-10-18"
# Updates the 2024 (column K) violent-crime counts for the Citywide Totals sheet,
# each affected neighborhood sheet, and the By Neighborhood summary sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 6363
$ws.Range("K3").Value = 6569
$ws.Range("K4").Value = 1367
$ws.Range("K6").Value = 7234
$ws.Range("K7").Value = 21997

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K2").Value = 191
$ws.Range("K5").Value = 57
$ws.Range("K6").Value = 157
$ws.Range("K7").Value = 647
$ws.Range("K8").Value = 1447
$ws.Range("K10").Value = 129
$ws.Range("K11").Value = 409
$ws.Range("K15").Value = 224
$ws.Range("K19").Value = 641
$ws.Range("K20").Value = 522
$ws.Range("K22").Value = 67
$ws.Range("K27").Value = 207
$ws.Range("K29").Value = 1189
$ws.Range("K30").Value = 83
$ws.Range("K31").Value = 245
$ws.Range("K33").Value = 961
$ws.Range("K37").Value = 748
$ws.Range("K41").Value = 154
$ws.Range("K42").Value = 814
$ws.Range("K43").Value = 180
$ws.Range("K44").Value = 182
$ws.Range("K46").Value = 44
$ws.Range("K48").Value = 275
$ws.Range("K51").Value = 282
$ws.Range("K52").Value = 581
$ws.Range("K54").Value = 432
$ws.Range("K63").Value = 59
$ws.Range("K67").Value = 862
$ws.Range("K72").Value = 113
$ws.Range("K76").Value = 300
$ws.Range("K77").Value = 151
$ws.Range("K78").Value = 250
$ws.Range("K79").Value = 554
$ws.Range("K85").Value = 1022
$ws.Range("K88").Value = 236
$ws.Range("K89").Value = 324
$ws.Range("K91").Value = 256
$ws.Range("K95").Value = 362
$ws.Range("K96").Value = 232
$ws.Range("K97").Value = 175
$ws.Range("K99").Value = 362
$ws.Range("K101").Value = 21997

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K6").Value = 98
$ws.Range("K7").Value = 232

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K2").Value = 211
$ws.Range("K6").Value = 178
$ws.Range("K7").Value = 647

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K6").Value = 134
$ws.Range("K7").Value = 409

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K3").Value = 99
$ws.Range("K7").Value = 324

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 335
$ws.Range("K3").Value = 352
$ws.Range("K7").Value = 1022

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K2").Value = 153
$ws.Range("K3").Value = 166
$ws.Range("K4").Value = 33
$ws.Range("K6").Value = 210
$ws.Range("K7").Value = 581

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 396
$ws.Range("K3").Value = 441
$ws.Range("K6").Value = 487
$ws.Range("K7").Value = 1447

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K4").Value = 50
$ws.Range("K6").Value = 296
$ws.Range("K7").Value = 961

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K2").Value = 123
$ws.Range("K7").Value = 362

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K2").Value = 211
$ws.Range("K3").Value = 248
$ws.Range("K6").Value = 223
$ws.Range("K7").Value = 748

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("K6").Value = 92
$ws.Range("K7").Value = 362

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("K3").Value = 26
$ws.Range("K7").Value = 83

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("K2").Value = 81
$ws.Range("K7").Value = 245

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K3").Value = 316
$ws.Range("K6").Value = 243
$ws.Range("K7").Value = 862

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K6").Value = 233
$ws.Range("K7").Value = 432

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K6").Value = 342
$ws.Range("K7").Value = 1189

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K4").Value = 39
$ws.Range("K6").Value = 129
$ws.Range("K7").Value = 275

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K3").Value = 193
$ws.Range("K7").Value = 641

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("K2").Value = 51
$ws.Range("K7").Value = 182

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K4").Value = 19
$ws.Range("K7").Value = 300

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("K6").Value = 43
$ws.Range("K7").Value = 157

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("K2").Value = 52
$ws.Range("K7").Value = 154

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K2").Value = 221
$ws.Range("K3").Value = 248
$ws.Range("K7").Value = 814

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("K6").Value = 58
$ws.Range("K7").Value = 129

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K2").Value = 75
$ws.Range("K4").Value = 23
$ws.Range("K7").Value = 250

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("K6").Value = 14
$ws.Range("K7").Value = 44

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K2").Value = 65
$ws.Range("K4").Value = 13
$ws.Range("K7").Value = 256

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K3").Value = 179
$ws.Range("K6").Value = 140
$ws.Range("K7").Value = 554

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K2").Value = 181
$ws.Range("K7").Value = 522

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("K2").Value = 82
$ws.Range("K6").Value = 68
$ws.Range("K7").Value = 224

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("K6").Value = 61
$ws.Range("K7").Value = 191

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("K2").Value = 36
$ws.Range("K3").Value = 37
$ws.Range("K7").Value = 175

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("K2").Value = 61
$ws.Range("K3").Value = 71
$ws.Range("K6").Value = 96
$ws.Range("K7").Value = 236

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("K2").Value = 14
$ws.Range("K7").Value = 57

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("K2").Value = 54
$ws.Range("K7").Value = 207

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("K2").Value = 78
$ws.Range("K7").Value = 282

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("K3").Value = 47
$ws.Range("K7").Value = 180

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("K3").Value = 20
$ws.Range("K7").Value = 67

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("K2").Value = 25
$ws.Range("K6").Value = 54
$ws.Range("K7").Value = 113

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("K4").Value = 11
$ws.Range("K7").Value = 151
